$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph
$range = $d.Content
$found = $range.Find.Execute("Docente(s) Responsável(eis)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Resolve which paragraph (by 1-based index) contains the found text. Using
# the paragraph's own Range (instead of the Find range, which can exclude a
# paragraph's trailing space/characters) keeps the heading text untouched.
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $range.Start -and $range.End -le $p.Range.End) {
        $headingIndex = $i
    }
}

$headingPara = $d.Paragraphs.Item($headingIndex)

# Insert a new paragraph mark right after the heading paragraph.
$insertRange = $headingPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

# The newly created paragraph is now immediately after the heading paragraph.
$newPara = $d.Paragraphs.Item($headingIndex + 1)
$newPara.Style = "List Bullet"
$newPara.Range.Text = "7455355 - Robson da Silva Rocha"
